$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Rename column O header from Num_Isolates to Present_SR
    $ws.Range("O1").Value = "Present_SR"

    # Delete the entire AO column (single_lineage); AP (Phenos) shifts left to become AO
    $result = $ws.Columns("AO").Delete()
}

Write-Output "done"
